$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p103r_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p103r_1</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p103r_2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p103r_2</id>", 2) | Out-Null
